# Update ref luong cua Phu tai CAN THO
# The "Phu cap" (allowance) column on the "Luong co ban" sheet is removed,
# shifting the location-percentage columns (CAN THO / LONG XUYEN / SOC TRANG)
# one column to the left.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Lương cơ bản"

# Remove the whole "Phụ cấp" column (R) from the main sheet; columns to the
# right (CẦN THƠ / LONG XUYÊN / SÓC TRĂNG) shift left automatically.
$ws1.Columns("R").Delete()

# Make "Lương cơ bản" the active sheet/tab with the given selection, which
# also clears the previous tab selection on "Chiết khấu".
$ws1.Activate()
$ws1.Range("T10").Select()
